$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "Datos actualizados..." timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 23 de Julio de 2020 a las 08:14"

# Row 6 - India
$ws.Range("B6").Value = 1239693
$ws.Range("C6").Value = 9
$ws.Range("E6").Value = 425537

# Row 32 - Kazajistan
$ws.Range("B32").Value = 76799
$ws.Range("C32").Value = 1646
$ws.Range("D32").Value = 48220
$ws.Range("E32").Value = 27994

# Row 53 - Afganistan
$ws.Range("B53").Value = 35915
$ws.Range("C53").Value = 188
$ws.Range("D53").Value = 24538
$ws.Range("E53").Value = 10166
$ws.Range("G53").Value = 21
$ws.Range("H53").Value = 1211

# Row 65 - Uzbekistan
$ws.Range("B65").Value = 18531
$ws.Range("C65").Value = 152
$ws.Range("E65").Value = 8560
$ws.Range("G65").Value = 1
$ws.Range("H65").Value = 99

# Row 76 - El Salvador
$ws.Range("D76").Value = 7151
$ws.Range("E76").Value = 5452
$ws.Range("G76").Value = 9
$ws.Range("H76").Value = 372

# Row 116 - Montenegro
$ws.Range("D116").Value = 496
$ws.Range("E116").Value = 1937
